$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# The chapter document was regenerated from source and the trailing
# "blog furniture" (social/category tags, post navigation, reader
# comments, "Leave a reply" form, Akismet notice, etc.) that used to
# follow the References list is no longer produced. What remains is the
# paragraph that used to hold a single trailing space right after the
# References list; it now carries a few extra blank runs where the
# removed content used to begin.
# -----------------------------------------------------------------------

# Locate the paragraph that immediately follows the References list
# item for "ValidationAttribute Class ..." - i.e. the lone-space
# paragraph - by searching for the last bullet in that list.
$markerRange = $d.Content.Duplicate
$found = $markerRange.Find.Execute(
    "ValidationAttribute Class (System.ComponentModel.DataAnnotations)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the References list anchor paragraph"
}

$anchorPara = $d.Paragraphs.Item($markerRange.Paragraphs.First.Index)
$keepPara = $anchorPara.Next()

# Everything from the paragraph right after $keepPara through to the end
# of the document is the legacy "blog furniture" being dropped.
$firstParaToDelete = $keepPara.Next()
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

if ($firstParaToDelete.Range.Start -le $lastPara.Range.End) {
    $deleteRange = $d.Range($firstParaToDelete.Range.Start, $lastPara.Range.End)
    $deleteRange.Delete()
}

# The surviving paragraph (formerly just a single space) now ends the
# document; give it three separate single-space runs (instead of the one
# it had) so it still renders as a small amount of trailing whitespace.
$target = $d.Range($keepPara.Range.Start, $keepPara.Range.End - 1)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlSnippet)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
Write-Output "Final paragraph text: [$($d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)]"
